$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.688.94"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.789.54"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.00"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4542"
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3697"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07305"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8539"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.36"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").Value = "1.800.92"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.560"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.302"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07067"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.80"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008626"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.64"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "26.706.39"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.262"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.71"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "2.026.53"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.71"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.171"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.26"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.204"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.59"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08839"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7597"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.154"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.447"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.887"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.118"
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01943"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05189"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.134"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.854"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.330"
$ws.Range("E42").Value = "  +6.88%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5210"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1651"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.463"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4953"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.17"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.77"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.648"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06281"
$ws.Range("E51").Value = "  -1.15%  "
